# natmiOut/YoungD7/LR-pairs_lrc2p/Efna4-Epha1.xlsx - "Natmi following Dr Hou advice"
#
# The sending/target cluster set changed from {FAPs, sCs} to {ECs, FAPs, sCs}
# (3 clusters x 4 target clusters = 12 data rows instead of 2 x 4 = 8), and every
# numeric column was recomputed. This rewrites rows 2-13 of Sheet1 in place,
# cell by cell, with the new values; the sheet/shared-string tables and the
# dimension ref are maintained automatically by Excel as the cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Efna4/Epha1)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna4"
$ws.Range("C2").Value = "Epha1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6731353333333333
$ws.Range("H2").Value = 2.019406
$ws.Range("I2").Value = 0.3272865828458516
$ws.Range("J2").Value = 0.3272865828458516
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.889554666666666
$ws.Range("N2").Value = 17.668664
$ws.Range("O2").Value = 0.295601531529132
$ws.Range("P2").Value = 0.295601531529132
$ws.Range("Q2").Value = 3.964467343731555
$ws.Range("R2").Value = 35.680206093584
$ws.Range("S2").Value = 0.09674641513816987
$ws.Range("T2").Value = 0.0967464151381699

# Row 3: ECs -> FAPs (Efna4/Epha1)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna4"
$ws.Range("C3").Value = "Epha1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.6731353333333333
$ws.Range("H3").Value = 2.019406
$ws.Range("I3").Value = 0.3272865828458516
$ws.Range("J3").Value = 0.3272865828458516
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.851268999999999
$ws.Range("N3").Value = 17.553807
$ws.Range("O3").Value = 0.2936799428279805
$ws.Range("P3").Value = 0.2936799428279805
$ws.Range("Q3").Value = 3.938695908738
$ws.Range("R3").Value = 35.448263178642
$ws.Range("S3").Value = 0.0961175049385348
$ws.Range("T3").Value = 0.09611750493853481

# Row 4: ECs -> M2 (Efna4/Epha1)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna4"
$ws.Range("C4").Value = "Epha1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6731353333333333
$ws.Range("H4").Value = 2.019406
$ws.Range("I4").Value = 0.3272865828458516
$ws.Range("J4").Value = 0.3272865828458516
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.039576333333334
$ws.Range("N4").Value = 6.118729000000001
$ws.Range("O4").Value = 0.1023679924759288
$ws.Range("P4").Value = 0.1023679924759288
$ws.Range("Q4").Value = 1.372910894997111
$ws.Range("R4").Value = 12.356198054974
$ws.Range("S4").Value = 0.0335036704502366
$ws.Range("T4").Value = 0.0335036704502366

# Row 5: ECs -> sCs (Efna4/Epha1)
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna4"
$ws.Range("C5").Value = "Epha1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.6731353333333333
$ws.Range("H5").Value = 2.019406
$ws.Range("I5").Value = 0.3272865828458516
$ws.Range("J5").Value = 0.3272865828458516
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.143565333333334
$ws.Range("N5").Value = 18.430696
$ws.Range("O5").Value = 0.3083505331669586
$ws.Range("P5").Value = 0.3083505331669587
$ws.Range("Q5").Value = 4.135450898508445
$ws.Range("R5").Value = 37.21905808657601
$ws.Range("S5").Value = 0.1009189923189103
$ws.Range("T5").Value = 0.1009189923189103

# Row 6: FAPs -> ECs (Efna4/Epha1)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna4"
$ws.Range("C6").Value = "Epha1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9964423333333334
$ws.Range("H6").Value = 2.989327
$ws.Range("I6").Value = 0.4844823769162027
$ws.Range("J6").Value = 0.4844823769162026
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.889554666666666
$ws.Range("N6").Value = 17.668664
$ws.Range("O6").Value = 0.295601531529132
$ws.Range("P6").Value = 0.295601531529132
$ws.Range("Q6").Value = 5.868601594347556
$ws.Range("R6").Value = 52.81741434912801
$ws.Range("S6").Value = 0.1432137326153037
$ws.Range("T6").Value = 0.1432137326153037

# Row 7: FAPs -> FAPs (Efna4/Epha1)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna4"
$ws.Range("C7").Value = "Epha1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9964423333333334
$ws.Range("H7").Value = 2.989327
$ws.Range("I7").Value = 0.4844823769162027
$ws.Range("J7").Value = 0.4844823769162026
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.851268999999999
$ws.Range("N7").Value = 17.553807
$ws.Range("O7").Value = 0.2936799428279805
$ws.Range("P7").Value = 0.2936799428279805
$ws.Range("Q7").Value = 5.830452135321
$ws.Range("R7").Value = 52.474069217889
$ws.Range("S7").Value = 0.1422827567539145
$ws.Range("T7").Value = 0.1422827567539145

# Row 8: FAPs -> M2 (Efna4/Epha1)
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna4"
$ws.Range("C8").Value = "Epha1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9964423333333334
$ws.Range("H8").Value = 2.989327
$ws.Range("I8").Value = 0.4844823769162027
$ws.Range("J8").Value = 0.4844823769162026
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.039576333333334
$ws.Range("N8").Value = 6.118729000000001
$ws.Range("O8").Value = 0.1023679924759288
$ws.Range("P8").Value = 0.1023679924759288
$ws.Range("Q8").Value = 2.032320200598112
$ws.Range("R8").Value = 18.29088180538301
$ws.Range("S8").Value = 0.04959548831487795
$ws.Range("T8").Value = 0.04959548831487795

# Row 9: FAPs -> sCs (Efna4/Epha1)
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna4"
$ws.Range("C9").Value = "Epha1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9964423333333334
$ws.Range("H9").Value = 2.989327
$ws.Range("I9").Value = 0.4844823769162027
$ws.Range("J9").Value = 0.4844823769162026
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.143565333333334
$ws.Range("N9").Value = 18.430696
$ws.Range("O9").Value = 0.3083505331669586
$ws.Range("P9").Value = 0.3083505331669587
$ws.Range("Q9").Value = 6.121708575732446
$ws.Range("R9").Value = 55.09537718159201
$ws.Range("S9").Value = 0.1493903992321065
$ws.Range("T9").Value = 0.1493903992321065

# Row 10: sCs -> ECs (Efna4/Epha1)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna4"
$ws.Range("C10").Value = "Epha1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3871376666666667
$ws.Range("H10").Value = 1.161413
$ws.Range("I10").Value = 0.1882310402379457
$ws.Range("J10").Value = 0.1882310402379457
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.889554666666666
$ws.Range("N10").Value = 17.668664
$ws.Range("O10").Value = 0.295601531529132
$ws.Range("P10").Value = 0.295601531529132
$ws.Range("Q10").Value = 2.280068451359111
$ws.Range("R10").Value = 20.520616062232
$ws.Range("S10").Value = 0.05564138377565843
$ws.Range("T10").Value = 0.05564138377565844

# Row 11: sCs -> FAPs (Efna4/Epha1)
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Efna4"
$ws.Range("C11").Value = "Epha1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3871376666666667
$ws.Range("H11").Value = 1.161413
$ws.Range("I11").Value = 0.1882310402379457
$ws.Range("J11").Value = 0.1882310402379457
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.851268999999999
$ws.Range("N11").Value = 17.553807
$ws.Range("O11").Value = 0.2936799428279805
$ws.Range("P11").Value = 0.2936799428279805
$ws.Range("Q11").Value = 2.265246627699
$ws.Range("R11").Value = 20.387219649291
$ws.Range("S11").Value = 0.0552796811355312
$ws.Range("T11").Value = 0.05527968113553121

# Row 12: sCs -> M2 (Efna4/Epha1)
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efna4"
$ws.Range("C12").Value = "Epha1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3871376666666667
$ws.Range("H12").Value = 1.161413
$ws.Range("I12").Value = 0.1882310402379457
$ws.Range("J12").Value = 0.1882310402379457
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.039576333333334
$ws.Range("N12").Value = 6.118729000000001
$ws.Range("O12").Value = 0.1023679924759288
$ws.Range("P12").Value = 0.1023679924759288
$ws.Range("Q12").Value = 0.7895968226752224
$ws.Range("R12").Value = 7.106371404077001
$ws.Range("S12").Value = 0.01926883371081429
$ws.Range("T12").Value = 0.01926883371081429

# Row 13: sCs -> sCs (Efna4/Epha1)
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efna4"
$ws.Range("C13").Value = "Epha1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3871376666666667
$ws.Range("H13").Value = 1.161413
$ws.Range("I13").Value = 0.1882310402379457
$ws.Range("J13").Value = 0.1882310402379457
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 6.143565333333334
$ws.Range("N13").Value = 18.430696
$ws.Range("O13").Value = 0.3083505331669586
$ws.Range("P13").Value = 0.3083505331669587
$ws.Range("Q13").Value = 2.378405548160889
$ws.Range("R13").Value = 21.405649933448
$ws.Range("S13").Value = 0.0580411416159418
$ws.Range("T13").Value = 0.05804114161594182
